# Updated symbol list on Sat Jan  7 15:37:16 UTC 2023 with GitHub Actions
# Refreshes the Price (col D) and Volume(1h) (col E) columns with the
# latest scraped quotes. Values are stored as plain text in the sheet
# (e.g. "260.76", "1.52%"), so each cell is written with a leading
# apostrophe to force text entry and the cell's original style is
# restored afterwards so no incidental number-format/style drift occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeRef, $value) {
    $cell = $ws.Range($rangeRef)
    $origStyle = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $origStyle
}

Set-TextValue "D2" "260.76"
Set-TextValue "E2" "1.52%"

Set-TextValue "D3" "27.16"
Set-TextValue "E3" "1.62%"

Set-TextValue "D4" "4.672"
Set-TextValue "E4" "0.73%"

Set-TextValue "E5" "4.22%"

Set-TextValue "D6" "6.673"
Set-TextValue "E6" "1.00%"

Set-TextValue "D7" "0.8514"
Set-TextValue "E7" "-0.61%"

Set-TextValue "D8" "0.9168"
Set-TextValue "E8" "-0.06%"

Set-TextValue "E9" "2.27%"

Set-TextValue "D10" "0.04801"
Set-TextValue "E10" "7.89%"

Set-TextValue "D12" "0.03119"
Set-TextValue "E12" "3.24%"

Set-TextValue "D13" "0.09047"
Set-TextValue "E13" "-0.69%"

Set-TextValue "D14" "0.001542"
Set-TextValue "E14" "0.84%"

Set-TextValue "D15" "0.0006186"
Set-TextValue "E15" "2.09%"

Set-TextValue "D16" "0.006122"
Set-TextValue "E16" "-1.26%"

Set-TextValue "E17" "-0.49%"

Set-TextValue "D18" "3.153"
Set-TextValue "E18" "0.63%"

Set-TextValue "E21" "0.18%"

Set-TextValue "E22" "5.76%"

Set-TextValue "D23" "0.04235"
Set-TextValue "E23" "0.52%"

Set-TextValue "E24" "0.11%"

Set-TextValue "E25" "-15.04%"

Set-TextValue "E26" "0.11%"

Set-TextValue "D27" "0.0001574"
Set-TextValue "E27" "-8.16%"

Set-TextValue "D40" "0.03876"
Set-TextValue "E40" "1.78%"

Set-TextValue "D42" "0.004093"
Set-TextValue "E42" "10.91%"

Set-TextValue "E43" "8.65%"

Set-TextValue "D44" "0.002196"
Set-TextValue "E44" "-9.54%"

Set-TextValue "D45" "0.00005153"
Set-TextValue "E45" "1.29%"

Set-TextValue "E46" "0.03%"

Set-TextValue "E47" "7.85%"

Set-TextValue "E49" "0.03%"

Set-TextValue "E50" "0.03%"
